$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The banner text no longer mentions editing the chart title - only the data.
$ws.Range("A1").Value = "Edit this spreadsheet to alter embedded chart data"

# Highlight the banner row: A1:E1 get a yellow fill (RGB 255,255,0 = 65535),
# F1 gets a plain theme "Background 1" / white fill so the banner reads as a
# single highlighted block.
$ws.Range("A1:F1").Interior.Color = 65535
$ws.Range("F1").Interior.ThemeColor = 2

# Leave the selection where the author ended up after the review fixes.
[void]$ws.Range("I3").Select()
